$d = $word.ActiveDocument

# 1. Merge "...finding success" + "." into a single run (no visible text change).
#    Avoid touching the straight apostrophe in "project's" to prevent autocorrect
#    from turning it into a curly quote.
$d.Content.Find.Execute("discover some trick for finding success", $true, $false, $false, $false, $false, `
    $true, 1, $false, "discover some trick for finding success", 2)

# 2. Merge "A" + "nalyze a database..." into a single run (no visible text change).
$d.Content.Find.Execute("Analyze a database of 4,000 past projects in order to uncover any hidden trends.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Analyze a database of 4,000 past projects in order to uncover any hidden trends.", 2)

# 3. Merge " rate" + ".  T" into a single run (no visible text change).
$d.Content.Find.Execute(" rate.  T", $true, $false, $false, $false, $false, `
    $true, 1, $false, " rate.  T", 2)

# 4. Merge "), and " + "theater (58%)" into a single run (no visible text change).
$d.Content.Find.Execute("), and theater (58%)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "), and theater (58%)", 2)

# 5. Add the missing word "Plays" before "Is the most popular sub-category."
$d.Content.Find.Execute("Is the most popular sub-category.", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Plays is the most popular sub-category.", 2)

# 6. "levels even" -> "levels out"
$d.Content.Find.Execute(" even and decline again until the month of December", $true, $false, $false, $false, $false, `
    $true, 1, $false, " out and decline again until the month of December", 2)

# 7. "An interesting observation, since" -> "Since"
$d.Content.Find.Execute("An interesting observation, since theater has the largest number of projects by far, if the theater category is ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Since theater has the largest number of projects by far, if the theater category is ", 2)

# 7b. Append new sentence after "...February has the highest number of success.  "
$d.Content.Find.Execute(", February has the highest number of success.  ", $true, $false, $false, $false, $false, `
    $true, 1, $false, ", February has the highest number of success.  The season also seem to affect success rate.", 2)

# 8. Merge "...outcome" + "' but my result shows that" into a single run (no visible text change).
$d.Content.Find.Execute("only a third have made it through the funding process with a positive outcome" + [char]8217, $true, $false, $false, $false, $false, `
    $true, 1, $false, "only a third have made it through the funding process with a positive outcome" + [char]8217, 2)
